# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with freshly scraped figures. Source values are plain text
# (scraped HTML), so numeric-looking Price cells are written with a
# leading apostrophe to force a text entry and keep formatting such as
# trailing zeros (e.g. "0.0920") intact instead of being normalised by
# Excel's automatic number detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $ws.Range($CellRef).Value = "'" + $Value
}

# Price column (D) - values that Excel would otherwise coerce to numbers.
Set-TextValue "D5"  "273.07"
Set-TextValue "D6"  "86.81"
Set-TextValue "D9"  "0.604"
Set-TextValue "D10" "44.86"
Set-TextValue "D11" "0.0920"
Set-TextValue "D12" "7.72"
Set-TextValue "D15" "14.88"
Set-TextValue "D17" "0.791"
Set-TextValue "D23" "232.38"
Set-TextValue "D24" "8.66"
Set-TextValue "D26" "2.55"
Set-TextValue "D27" "10.77"
Set-TextValue "D30" "39.24"
Set-TextValue "D31" "172.76"
Set-TextValue "D32" "0.0904"
Set-TextValue "D33" "20.73"
Set-TextValue "D37" "0.0352"
Set-TextValue "D39" "3.36"
Set-TextValue "D41" "12.42"
Set-TextValue "D42" "63.57"
Set-TextValue "D43" "0.206"
Set-TextValue "D45" "8.48"
Set-TextValue "D47" "99.92"
Set-TextValue "D51" "1.48"

# Price column (D) - values that already contain two decimal points and
# therefore stay text automatically, but are set the same way for
# consistency.
Set-TextValue "D2"  "43.780.41"
Set-TextValue "D3"  "2.229.29"
Set-TextValue "D14" "2.564.14"
Set-TextValue "D16" "2.223.71"
Set-TextValue "D18" "43.713.62"

# Volume(1h) column (E) - percentage strings, never numeric, safe as-is.
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +5.82%  "
$ws.Range("E6").Value = "  +9.85%  "
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  +8.98%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  -8.35%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +14.04%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  +4.10%  "
$ws.Range("E29").Value = "  +5.33%  "
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("E39").Value = "  +16.57%  "
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("E41").Value = "  -5.94%  "
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -4.68%  "
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  -8.41%  "
$ws.Range("E51").Value = "  -1.82%  "
